# Insert the newest day's two rows (Primera / Segunda quality) for
# "Coliflor" at the top of the data block (row 1043), pushing the
# existing 1043:1170 block down to 1045:1172.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1043:1044").Insert()

# --- Row 1043: Primera ---
$ws.Range("A1043").Value = 8
$ws.Range("B1043").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1043").Value = "Coquimbo"
$ws.Range("D1043").Value = 45142
$ws.Range("E1043").Value = 4
$ws.Range("F1043").Value = 100112008
$ws.Range("G1043").Value = "Coliflor"
$ws.Range("H1043").Value = "Sin especificar"
$ws.Range("I1043").Value = "Primera"
$ws.Range("J1043").Value = 2000
$ws.Range("K1043").Value = 700
$ws.Range("L1043").Value = 800
$ws.Range("M1043").Value = 750
$ws.Range("N1043").Value = "$/unidad"
$ws.Range("O1043").Value = "Provincia del Elquí"
$ws.Range("P1043").Value = 750
$ws.Range("Q1043").Value = 1
$ws.Range("R1043").Value = "Hortaliza"

# --- Row 1044: Segunda ---
$ws.Range("A1044").Value = 8
$ws.Range("B1044").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1044").Value = "Coquimbo"
$ws.Range("D1044").Value = 45142
$ws.Range("E1044").Value = 4
$ws.Range("F1044").Value = 100112008
$ws.Range("G1044").Value = "Coliflor"
$ws.Range("H1044").Value = "Sin especificar"
$ws.Range("I1044").Value = "Segunda"
$ws.Range("J1044").Value = 1200
$ws.Range("K1044").Value = 500
$ws.Range("L1044").Value = 600
$ws.Range("M1044").Value = 550
$ws.Range("N1044").Value = "$/unidad"
$ws.Range("O1044").Value = "Provincia del Elquí"
$ws.Range("P1044").Value = 550
$ws.Range("Q1044").Value = 1
$ws.Range("R1044").Value = "Hortaliza"
